$d = $word.ActiveDocument

$pairs = @(
    @("3+94=", "85-35="),
    @("64+22=", "19+39="),
    @("15+17=", "46+53="),
    @("67-44=", "8+41="),
    @("18+59=", "18+58="),
    @("87-21=", "22+29="),
    @("61-32=", "55+21="),
    @("43+43=", "67-15="),
    @("55+35=", "24-3="),
    @("71-48=", "22+7="),
    @("63+5=", "16+28="),
    @("17-13=", "5+22="),
    @("10+1=", "31+68="),
    @("9+58=", "33-20="),
    @("64-21=", "40+13="),
    @("66+19=", "69-68="),
    @("0+26=", "90-80="),
    @("82+9=", "38-13="),
    @("6+66=", "99-89="),
    @("75-29=", "49+12="),
    @("24+63=", "53+29="),
    @("87-54=", "79+7="),
    @("6+4=", "4+17="),
    @("56+29=", "65-31="),
    @("51+28=", "18+55="),
    @("71+1=", "30+9="),
    @("13+24=", "37-15="),
    @("44-30=", "9+17="),
    @("20+51=", "57+28="),
    @("67+6=", "84-11="),
    @("8+31=", "57-55="),
    @("22+43=", "80-41="),
    @("56-47=", "12+49="),
    @("99-73=", "1+74="),
    @("62-48=", "31-13="),
    @("26+20=", "58+36="),
    @("57-11=", "15+21="),
    @("52+42=", "81-59="),
    @("34+21=", "88-5="),
    @("66-25=", "8+84="),
    @("11+22=", "42+22="),
    @("77-63=", "29-12="),
    @("93-16=", "93-60="),
    @("31+29=", "10+33="),
    @("94-57=", "64+25="),
    @("9-1=", "21+75="),
    @("25+65=", "81+13="),
    @("92-4=", "62-25="),
    @("74-69=", "38-38="),
    @("25+3=", "97-74="),
    @("45+39=", "72+5="),
    @("84-12=", "77+19="),
    @("82-48=", "52-34="),
    @("62-0=", "91+7="),
    @("26+21=", "96-64="),
    @("67+25=", "20+72="),
    @("91-31=", "1+24="),
    @("32+37=", "20+64="),
    @("78-36=", "81-20="),
    @("87-9=", "36-25="),
    @("25+43=", "98-54="),
    @("82-61=", "82-9="),
    @("31+20=", "52-26="),
    @("60+13=", "42-11="),
    @("15+74=", "83-53="),
    @("97-66=", "94-72="),
    @("22+3=", "84+4="),
    @("83-42=", "34+46="),
    @("90-45=", "6+30="),
    @("40-19=", "74+6="),
    @("43-40=", "16+23="),
    @("41-18=", "13+46="),
    @("76+6=", "46+22="),
    @("35+59=", "91-56="),
    @("21+63=", "67-40="),
    @("17+55=", "8+20="),
    @("14-1=", "96-29="),
    @("1+60=", "70-21="),
    @("37-33=", "47-29="),
    @("14+25=", "27-7="),
    @("25+61=", "10+5="),
    @("28-3=", "96-76="),
    @("94-13=", "74+9="),
    @("92-54=", "46+48="),
    @("22-1=", "39-17="),
    @("36-16=", "97-55="),
    @("2+58=", "31+51="),
    @("25+19=", "72-17="),
    @("30-28=", "93-71="),
    @("64-11=", "0+39="),
    @("17+54=", "67-48="),
    @("15+9=", "27+60="),
    @("71-25=", "88-71="),
    @("96-15=", "48-18="),
    @("1+30=", "4+54="),
    @("31-24=", "84-2="),
    @("4+41=", "26+24="),
    @("89-40=", "91+7="),
    @("7+46=", "65+20="),
    @("87-70=", "63-54=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
